$d = $word.ActiveDocument

# Each cell in the multiplication-facts table holds a unique "A\xD7B=C" string.
# Replace each old expression with its new counterpart using Find/Replace
# scoped to the whole document story range.
$pairs = @(
    ,@("87×82=7134", "97×62=6014")
    ,@("55×13=715", "41×79=3239")
    ,@("76×89=6764", "97×53=5141")
    ,@("85×13=1105", "43×80=3440")
    ,@("72×10=720", "50×40=2000")
    ,@("13×38=494", "60×72=4320")
    ,@("20×99=1980", "27×21=567")
    ,@("69×61=4209", "46×58=2668")
    ,@("43×81=3483", "53×92=4876")
    ,@("51×29=1479", "21×35=735")
    ,@("69×94=6486", "84×72=6048")
    ,@("91×96=8736", "90×60=5400")
    ,@("48×36=1728", "91×61=5551")
    ,@("70×14=980", "40×98=3920")
    ,@("33×99=3267", "65×83=5395")
    ,@("72×93=6696", "18×93=1674")
    ,@("90×77=6930", "67×89=5963")
    ,@("57×92=5244", "74×50=3700")
    ,@("73×19=1387", "11×59=649")
    ,@("17×88=1496", "74×99=7326")
    ,@("88×97=8536", "10×18=180")
    ,@("60×33=1980", "85×61=5185")
    ,@("75×54=4050", "56×91=5096")
    ,@("65×15=975", "66×97=6402")
    ,@("10×13=130", "18×43=774")
    ,@("87×61=5307", "42×87=3654")
    ,@("54×25=1350", "57×60=3420")
    ,@("17×18=306", "28×96=2688")
    ,@("95×81=7695", "99×62=6138")
    ,@("95×27=2565", "16×93=1488")
    ,@("87×57=4959", "58×84=4872")
    ,@("74×71=5254", "18×94=1692")
    ,@("30×11=330", "71×92=6532")
    ,@("42×94=3948", "66×93=6138")
    ,@("85×90=7650", "55×56=3080")
    ,@("24×52=1248", "79×11=869")
    ,@("90×74=6660", "63×84=5292")
    ,@("12×69=828", "28×39=1092")
    ,@("54×48=2592", "74×92=6808")
    ,@("33×54=1782", "42×98=4116")
    ,@("82×41=3362", "21×30=630")
    ,@("90×16=1440", "74×44=3256")
    ,@("95×59=5605", "18×67=1206")
    ,@("78×62=4836", "71×98=6958")
    ,@("73×76=5548", "62×68=4216")
    ,@("25×50=1250", "67×72=4824")
    ,@("29×41=1189", "86×29=2494")
    ,@("57×67=3819", "77×12=924")
    ,@("57×72=4104", "78×40=3120")
    ,@("49×61=2989", "76×46=3496")
    ,@("79×36=2844", "13×26=338")
    ,@("98×72=7056", "93×65=6045")
    ,@("38×46=1748", "51×52=2652")
    ,@("14×86=1204", "35×94=3290")
    ,@("66×17=1122", "81×16=1296")
    ,@("32×33=1056", "48×64=3072")
    ,@("99×29=2871", "39×37=1443")
    ,@("65×17=1105", "56×95=5320")
    ,@("94×40=3760", "17×19=323")
    ,@("49×31=1519", "25×35=875")
    ,@("67×12=804", "86×42=3612")
    ,@("54×87=4698", "92×52=4784")
    ,@("29×70=2030", "22×71=1562")
    ,@("58×46=2668", "55×60=3300")
    ,@("11×66=726", "51×87=4437")
    ,@("47×24=1128", "23×60=1380")
    ,@("45×86=3870", "34×86=2924")
    ,@("98×62=6076", "55×33=1815")
    ,@("40×52=2080", "74×77=5698")
    ,@("99×42=4158", "84×91=7644")
    ,@("91×40=3640", "85×16=1360")
    ,@("91×58=5278", "89×75=6675")
    ,@("64×58=3712", "66×77=5082")
    ,@("93×38=3534", "33×44=1452")
    ,@("89×63=5607", "25×45=1125")
    ,@("14×77=1078", "40×38=1520")
    ,@("57×81=4617", "45×44=1980")
    ,@("43×17=731", "62×56=3472")
    ,@("67×96=6432", "34×41=1394")
    ,@("30×53=1590", "58×76=4408")
    ,@("95×58=5510", "53×14=742")
    ,@("48×48=2304", "26×27=702")
    ,@("69×73=5037", "88×63=5544")
    ,@("34×88=2992", "76×49=3724")
    ,@("88×53=4664", "61×58=3538")
    ,@("47×88=4136", "98×55=5390")
    ,@("59×40=2360", "77×52=4004")
    ,@("72×68=4896", "60×51=3060")
    ,@("94×28=2632", "25×26=650")
    ,@("10×72=720", "46×92=4232")
    ,@("21×75=1575", "94×52=4888")
    ,@("88×81=7128", "15×69=1035")
    ,@("53×57=3021", "51×77=3927")
    ,@("72×95=6840", "94×42=3948")
    ,@("95×40=3800", "15×22=330")
    ,@("76×82=6232", "20×71=1420")
    ,@("33×64=2112", "44×70=3080")
    ,@("52×55=2860", "100×22=2200")
    ,@("62×49=3038", "46×94=4324")
    ,@("51×67=3417", "66×86=5676")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

$d.Save()
